$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.834.07"
$ws.Range("E2").Value = "  +0.83%  "
$ws.Range("D3").Value = "2.509.08"
$ws.Range("E3").Value = "  +0.80%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.54"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.10"
$ws.Range("E6").Value = "  +3.78%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +0.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.55"
$ws.Range("E10").Value = "  +6.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0814"
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.65"
$ws.Range("E13").Value = "  +1.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.26"
$ws.Range("E14").Value = "  +1.20%  "
$ws.Range("D15").Value = "2.902.89"
$ws.Range("E15").Value = "  +0.80%  "
$ws.Range("D16").Value = "2.505.75"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.852"
$ws.Range("E17").Value = "  +0.62%  "
$ws.Range("D18").Value = "47.819.23"
$ws.Range("E18").Value = "  +1.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.33"
$ws.Range("E19").Value = "  +4.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.65"
$ws.Range("E20").Value = "  +1.28%  "
$ws.Range("B21").Value = "ImmutableX"
$ws.Range("C21").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.80"
$ws.Range("E21").Value = "  +16.15%  "
$ws.Range("B22").Value = "ShibaInu"
$ws.Range("C22").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D22").Value = "0.0₃0945"
$ws.Range("E22").Value = "  +0.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.77"
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "247.93"
$ws.Range("E24").Value = "  -1.30%  "
$ws.Range("E25").Value = "  -1.07%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.92"
$ws.Range("E27").Value = "  -1.01%  "
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("E30").Value = "  +3.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.07"
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.73"
$ws.Range("E32").Value = "  +0.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.18"
$ws.Range("E33").Value = "  +2.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.38"
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0789"
$ws.Range("E35").Value = "  +0.38%  "
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.70"
$ws.Range("E38").Value = "  +1.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.00"
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.71"
$ws.Range("E40").Value = "  +6.25%  "
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.112"
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("E42").Value = "  -1.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "119.18"
$ws.Range("E43").Value = "  -2.11%  "
$ws.Range("E44").Value = "  +0.57%  "
$ws.Range("D45").Value = "1.999.80"
$ws.Range("E45").Value = "  +1.78%  "
$ws.Range("E46").Value = "  +2.91%  "
$ws.Range("E47").Value = "  -3.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.84"
$ws.Range("E48").Value = "  +1.94%  "
$ws.Range("E49").Value = "  -1.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.23"
$ws.Range("E50").Value = "  -0.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "56.82"
$ws.Range("E51").Value = "  +3.76%  "
